$wb = $excel.ActiveWorkbook

# Sheet ALC, row 4 (Leve Item ID 5470)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# Sheet ALC, row 17 (Leve Item ID 38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1066.5385
$ws.Range("I17").Value = 433.33334
$ws.Range("J17").Value = 1256.5
$ws.Range("K17").Value = 1300.00002
$ws.Range("L17").Value = 3769.5
$ws.Range("M17").Value = -1132.00002
$ws.Range("N17").Value = -4105.5

# Sheet ALC, row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6988.5835
$ws.Range("I76").Value = 9841.143
$ws.Range("J76").Value = 2995
$ws.Range("K76").Value = 9841.143
$ws.Range("L76").Value = 2995
$ws.Range("M76").Value = -9526.143

# Sheet ALC, row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6988.5835
$ws.Range("I79").Value = 9841.143
$ws.Range("J79").Value = 2995
$ws.Range("K79").Value = 9841.143
$ws.Range("L79").Value = 2995
$ws.Range("M79").Value = -8749.143

# Sheet ALC, row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1478.4
$ws.Range("I100").Value = 598.3333
$ws.Range("J100").Value = 1855.5714
$ws.Range("K100").Value = 598.3333
$ws.Range("L100").Value = 1855.5714
$ws.Range("M100").Value = -57.33330000000001

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4165.2144
$ws.Range("I132").Value = 4216.385
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 12649.155
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -10119.155
$ws.Range("N132").Value = -15560

# Sheet ALC, row 133 (Leve Item ID 41856)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 100780
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100780
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 100780
$ws.Range("N133").Value = -110900

# Sheet ALC, row 134 (Leve Item ID 41997)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 69780
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 69780
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 69780
$ws.Range("N134").Value = -79920

# Sheet ALC, row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4599.615
$ws.Range("I138").Value = 900.8889
$ws.Range("J138").Value = 8594.24
$ws.Range("K138").Value = 2702.6667
$ws.Range("L138").Value = 25782.72
$ws.Range("M138").Value = 2437.3333
$ws.Range("N138").Value = -36062.72

# Sheet ALC, row 141 (Leve Item ID 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1675.125
$ws.Range("I141").Value = 1480.3158
$ws.Range("J141").Value = 2415.4
$ws.Range("K141").Value = 4440.9474
$ws.Range("L141").Value = 7246.200000000001
$ws.Range("M141").Value = 739.0526

# Sheet ARM, row 23 (Leve Item ID 2236)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 25000
$ws.Range("I23").Value = 25000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -24741

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2554421.5
$ws.Range("I32").Value = 2720601.8
$ws.Range("J32").Value = 6322
$ws.Range("K32").Value = 2720601.8
$ws.Range("L32").Value = 6322
$ws.Range("M32").Value = -2720314.8

# Sheet ARM, row 43 (Leve Item ID 21715)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 14990
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 14990
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 14990
$ws.Range("N43").Value = -15616

# Sheet ARM, row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3585.9524
$ws.Range("I61").Value = 2649.8572
$ws.Range("J61").Value = 5458.143
$ws.Range("K61").Value = 2649.8572
$ws.Range("L61").Value = 5458.143
$ws.Range("M61").Value = -2437.8572
$ws.Range("N61").Value = -5882.143

# Sheet ARM, row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 33589.72
$ws.Range("I74").Value = 39101.047
$ws.Range("J74").Value = 4655.25
$ws.Range("K74").Value = 39101.047
$ws.Range("L74").Value = 4655.25
$ws.Range("M74").Value = -38227.047
$ws.Range("N74").Value = -6403.25

# Sheet ARM, row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 33589.72
$ws.Range("I77").Value = 39101.047
$ws.Range("J77").Value = 4655.25
$ws.Range("K77").Value = 195505.235
$ws.Range("L77").Value = 23276.25
$ws.Range("M77").Value = -191137.235
$ws.Range("N77").Value = -32012.25

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9997.5
$ws.Range("I132").Value = 10346.846
$ws.Range("J132").Value = 9694.733
$ws.Range("K132").Value = 31040.538
$ws.Range("L132").Value = 29084.199
$ws.Range("M132").Value = -28510.538

# Sheet ARM, row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3585.9524
$ws.Range("I136").Value = 2649.8572
$ws.Range("J136").Value = 5458.143
$ws.Range("K136").Value = 7949.571599999999
$ws.Range("L136").Value = 16374.429
$ws.Range("M136").Value = -5399.571599999999
$ws.Range("N136").Value = -21474.429

# Sheet BSM, row 99 (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 12989378

# Sheet BSM, row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4040.5134
$ws.Range("I134").Value = 2040.772
$ws.Range("J134").Value = 10745.529
$ws.Range("K134").Value = 6122.316
$ws.Range("L134").Value = 32236.587
$ws.Range("M134").Value = -3587.316

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6317.0625
$ws.Range("I31").Value = 2692.375
$ws.Range("J31").Value = 12358.208
$ws.Range("K31").Value = 2692.375
$ws.Range("L31").Value = 12358.208
$ws.Range("M31").Value = -2397.375

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6317.0625
$ws.Range("I34").Value = 2692.375
$ws.Range("J34").Value = 12358.208
$ws.Range("K34").Value = 2692.375
$ws.Range("L34").Value = 12358.208
$ws.Range("M34").Value = -2490.375

# Sheet CRP, row 37 (Leve Item ID 2021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 5100
$ws.Range("I37").Value = 5100
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5100
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -4993

# Sheet CRP, row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13519696
$ws.Range("I58").Value = 29413692
$ws.Range("J58").Value = 9800.15
$ws.Range("K58").Value = 29413692
$ws.Range("L58").Value = 9800.15
$ws.Range("M58").Value = -29413489
$ws.Range("N58").Value = -10206.15

# Sheet CRP, row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9569.200000000001
$ws.Range("I132").Value = 2548
$ws.Range("J132").Value = 14250
$ws.Range("K132").Value = 7644
$ws.Range("L132").Value = 42750
$ws.Range("M132").Value = -5114
$ws.Range("N132").Value = -47810

# Sheet CRP, row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8871.244000000001
$ws.Range("I134").Value = 10187.883
$ws.Range("J134").Value = 7938.625
$ws.Range("K134").Value = 30563.649
$ws.Range("L134").Value = 23815.875
$ws.Range("M134").Value = -28028.649
$ws.Range("N134").Value = -28885.875

# Sheet CRP, row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13519696
$ws.Range("I136").Value = 29413692
$ws.Range("J136").Value = 9800.15
$ws.Range("K136").Value = 88241076
$ws.Range("L136").Value = 29400.45
$ws.Range("M136").Value = -88238526
$ws.Range("N136").Value = -34500.45

# Sheet CUL, row 15 (Leve Item ID 4661)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 983.7143
$ws.Range("I15").Value = 95
$ws.Range("J15").Value = 1339.2
$ws.Range("K15").Value = 285
$ws.Range("L15").Value = 4017.6
$ws.Range("M15").Value = -145
$ws.Range("N15").Value = -4297.6

# Sheet CUL, row 26 (Leve Item ID 4746)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 463
$ws.Range("I26").Value = 49.5
$ws.Range("J26").Value = 554.8889
$ws.Range("K26").Value = 148.5
$ws.Range("L26").Value = 1664.6667
$ws.Range("M26").Value = 139.5
$ws.Range("N26").Value = -2240.6667

# Sheet GSM, row 21 (Leve Item ID 4430)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 7500
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -4827
$ws.Range("N21").Value = -10346

# Sheet GSM, row 30 (Leve Item ID 4430)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 7500
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = -4895
$ws.Range("N30").Value = -10210

# Sheet GSM, row 80 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 95109.55
$ws.Range("I80").Value = 4371.8335
$ws.Range("J80").Value = 203994.8
$ws.Range("K80").Value = 4371.8335
$ws.Range("L80").Value = 203994.8
$ws.Range("M80").Value = -3373.8335
$ws.Range("N80").Value = -205990.8

# Sheet GSM, row 83 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 95109.55
$ws.Range("I83").Value = 4371.8335
$ws.Range("J83").Value = 203994.8
$ws.Range("K83").Value = 21859.1675
$ws.Range("L83").Value = 1019974
$ws.Range("M83").Value = -16867.1675
$ws.Range("N83").Value = -1029958

# Sheet GSM, row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 47495.086
$ws.Range("I122").Value = 62011.65
$ws.Range("J122").Value = 6364.8335
$ws.Range("K122").Value = 186034.95
$ws.Range("L122").Value = 19094.5005
$ws.Range("M122").Value = -183584.95

# Sheet GSM, row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2572.1538
$ws.Range("I126").Value = 2701.3333
$ws.Range("J126").Value = 2533.4
$ws.Range("K126").Value = 8103.999899999999
$ws.Range("L126").Value = 7600.200000000001
$ws.Range("M126").Value = -5633.999899999999
$ws.Range("N126").Value = -12540.2

# Sheet GSM, row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2057.9148
$ws.Range("I132").Value = 2065.2
$ws.Range("J132").Value = 2016.2858
$ws.Range("K132").Value = 6195.599999999999
$ws.Range("L132").Value = 6048.857400000001
$ws.Range("M132").Value = -3665.599999999999
$ws.Range("N132").Value = -11108.8574

# Sheet LTW, row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5401.7
$ws.Range("I40").Value = 4895.0713
$ws.Range("J40").Value = 6583.8335
$ws.Range("K40").Value = 4895.0713
$ws.Range("L40").Value = 6583.8335
$ws.Range("M40").Value = -4759.0713
$ws.Range("N40").Value = -6855.8335

# Sheet LTW, row 68 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2714.8572
$ws.Range("I68").Value = 2400.6667
$ws.Range("J68").Value = 2950.5
$ws.Range("K68").Value = 2400.6667
$ws.Range("L68").Value = 2950.5
$ws.Range("M68").Value = -1651.6667
$ws.Range("N68").Value = -4448.5

# Sheet LTW, row 71 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2714.8572
$ws.Range("I71").Value = 2400.6667
$ws.Range("J71").Value = 2950.5
$ws.Range("K71").Value = 12003.3335
$ws.Range("L71").Value = 14752.5
$ws.Range("M71").Value = -8259.333500000001
$ws.Range("N71").Value = -22240.5

# Sheet LTW, row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3930.82
$ws.Range("I122").Value = 3687.5676
$ws.Range("J122").Value = 4623.154
$ws.Range("K122").Value = 11062.7028
$ws.Range("L122").Value = 13869.462
$ws.Range("M122").Value = -8612.702799999999

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20841444
$ws.Range("I132").Value = 45459820
$ws.Range("J132").Value = 10513.154
$ws.Range("K132").Value = 136379460
$ws.Range("L132").Value = 31539.462
$ws.Range("M132").Value = -136376930

# Sheet LTW, row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8337.647999999999
$ws.Range("I136").Value = 4074.4
$ws.Range("J136").Value = 13353.235
$ws.Range("K136").Value = 12223.2
$ws.Range("L136").Value = 40059.705
$ws.Range("M136").Value = -9673.200000000001

# Sheet WVR, row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3121.963
$ws.Range("I126").Value = 1464.65
$ws.Range("J126").Value = 7857.143
$ws.Range("K126").Value = 4393.950000000001
$ws.Range("L126").Value = 23571.429
$ws.Range("M126").Value = -1923.950000000001
$ws.Range("N126").Value = -28511.429

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 55625876
$ws.Range("I132").Value = 100024890
$ws.Range("J132").Value = 127111
$ws.Range("K132").Value = 300074670
$ws.Range("L132").Value = 381333
$ws.Range("M132").Value = -300072140
$ws.Range("N132").Value = -386393

# Sheet WVR, row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 41670690
$ws.Range("I136").Value = 76923976
$ws.Range("J136").Value = 7706.364
$ws.Range("K136").Value = 230771928
$ws.Range("L136").Value = 23119.092
$ws.Range("M136").Value = -230769378
